$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "category"
$ws.Range("B1").Value = "Amount"
$ws.Range("C1").Value = "Date"

# Data rows. The Date column holds DD/MM/YYYY-looking strings that must be
# kept as literal text (not auto-converted to real Excel date serials), so
# the cell is temporarily switched to Text format while the value is
# entered and then restored to the default "Normal" style so no stray
# number formatting is left behind.
$ws.Range("A2").Value = "Travel"
$ws.Range("B2").Value = 5000
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "11/08/2025"
$ws.Range("C2").Style = "Normal"

$ws.Range("A3").Value = "Rent"
$ws.Range("B3").Value = 20000
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "04/08/2025"
$ws.Range("C3").Style = "Normal"

$ws.Range("A4").Value = "tickets"
$ws.Range("B4").Value = 8000
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "02/08/2025"
$ws.Range("C4").Style = "Normal"
